# calculateRewards update and test ok
#
# Applies the "TourGuide Performance Graphs" edit:
#   - Locations sheet: row 4 height tweak, new D9 value, cursor moved to D4
#   - Rewards sheet:   D3:D6 reward counts recalculated (bigger numbers),
#                       new E6 / E10 cells, column C widened, cursor moved to D4
#
# (Chart axis-id bookkeeping and drawing-frame pixel offsets are internal
#  Excel-generated ids/caches that get regenerated from the sheet data the
#  next time the chart is opened/rendered - they are not reachable through
#  the Excel object model itself, so this script concentrates on the
#  underlying cell data and sheet/view properties that actually drive them.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locations sheet
# ---------------------------------------------------------------------
$locations = $wb.Worksheets.Item("Locations")

# Row 4 gets a slightly tighter custom height.
$locations.Rows.Item(4).RowHeight = 13.8

# New row: D9 = 900 (extends the used range to A1:D9).
$locations.Range("D9").Value = 900

# Move the cursor to D4 on this sheet.
$locations.Range("D4").Select()

# ---------------------------------------------------------------------
# Rewards sheet
# ---------------------------------------------------------------------
$rewards = $wb.Worksheets.Item("Rewards")

# Recalculated reward counts (D3:D6).
$rewards.Range("D3").Value = 3
$rewards.Range("D4").Value = 11
$rewards.Range("D5").Value = 2600
$rewards.Range("D6").Value = 26000

# New cells: E6 (a number) and E10 (a short text label, "unt" - as in
# "amoUNT"/"coUNT" - pulled in as a new shared string).
$rewards.Range("E6").Value = 1200
$rewards.Range("E10").Value = "unt"

# Column C is widened a bit for the new data.
$rewards.Columns.Item(3).ColumnWidth = 16.6

# Move the cursor to D4 on this sheet too, and make sure Rewards - which
# was already the active tab before the edit - stays the active tab.
$rewards.Range("D4").Select()
$rewards.Activate()
